$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: date (copy formatting from an existing date cell, then set the value)
$ws.Range("A27").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value2 = 43875

# Row 31/32: task text
$ws.Range("B31").Value = "Refactor write and read zarr file code"
$ws.Range("B32").Value = "Plan for incorporating napari pipeline with napari movie"

# Match the final view/selection state from the diff
$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 9
